$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title heading.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play 3 Stars Slot Game Free - Review & Demo"
#    right before the final paragraph (the one that currently holds the
#    feature-image prompt text).
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)
$insertPos = $finalPara.Range.Start - 1

$snippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 3 Stars Slot Game Free - Review &amp; Demo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertRng = $d.Range($insertPos, $insertPos)
$insertRng.InsertXML($snippet) | Out-Null

# 3. Replace the final paragraph's text (the old image-prompt text) with the
#    meta-description sentence, keeping its existing italic run formatting.
$d.Content.Find.Execute(
    "Create a feature image*development company.",
    $false, $false, $true, $false, $false, $true, 1, $false,
    "Take a spin and play 3 Stars, a 5-reel video slot game with 50 fixed paylines and four jackpots available to win. Review and demo available.",
    2) | Out-Null
